$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.035664916038513
$ws.Range("B1").Value = 0.9053687453269958
$ws.Range("C1").Value = 6.534774780273438
$ws.Range("D1").Value = 2.014665603637695
$ws.Range("E1").Value = 1.119258522987366
